$wb = $excel.ActiveWorkbook

# Switch the active sheet to "Repayment Schedule" (was "Transactions")
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank) column before column N, shifting old N/O/P -> O/P/Q
$ws.Columns("N").Insert()

# Update the selection on the Repayment Schedule sheet
$ws.Range("K19").Select()
